$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while forcing it to remain a text string,
# matching the original inlineStr/shared-string cell type, then restore
# the cell style so no stray formatting is introduced.
function Set-TextValue($address, $value) {
    $cell = $ws.Range($address)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextValue 'D2' '63.846.68'
Set-TextValue 'E2' '  +1.47%  '
Set-TextValue 'D3' '2.522.20'
Set-TextValue 'E3' '  +2.34%  '
Set-TextValue 'E4' '  +0.06%  '
Set-TextValue 'D5' '579.40'
Set-TextValue 'E5' '  +0.74%  '
Set-TextValue 'D6' '151.11'
Set-TextValue 'E6' '  +2.87%  '
Set-TextValue 'E7' '  +0.03%  '
Set-TextValue 'D8' '0.536'
Set-TextValue 'E8' '  -0.12%  '
Set-TextValue 'E9' '  +0.05%  '
Set-TextValue 'E10' '  -0.89%  '
Set-TextValue 'E11' '  -0.24%  '
Set-TextValue 'D12' '0.352'
Set-TextValue 'E12' '  -1.81%  '
Set-TextValue 'D13' '29.41'
Set-TextValue 'E13' '  +1.17%  '
Set-TextValue 'E14' '  -0.15%  '
Set-TextValue 'D15' '2.978.14'
Set-TextValue 'E15' '  +2.29%  '
Set-TextValue 'D16' '63.737.19'
Set-TextValue 'E16' '  +1.49%  '
Set-TextValue 'D17' '2.528.86'
Set-TextValue 'E17' '  +2.62%  '
Set-TextValue 'D18' '7.78'
Set-TextValue 'E18' '  -2.72%  '
Set-TextValue 'D19' '10.87'
Set-TextValue 'E19' '  -1.51%  '
Set-TextValue 'E20' '  +2.30%  '
Set-TextValue 'D21' '326.00'
Set-TextValue 'E21' '  -0.47%  '
Set-TextValue 'E22' '  +0.48%  '
Set-TextValue 'E23' '  -0.04%  '
Set-TextValue 'D24' '10.08'
Set-TextValue 'E24' '  -0.83%  '
Set-TextValue 'D25' '65.28'
Set-TextValue 'E25' '  -0.78%  '
Set-TextValue 'D26' '657.78'
Set-TextValue 'E26' '  +0.89%  '
Set-TextValue 'E27' '  +2.95%  '
Set-TextValue 'D28' '2.650.18'
Set-TextValue 'E28' '  +2.47%  '
Set-TextValue 'D29' '0.999'
Set-TextValue 'E29' '  +0.24%  '
Set-TextValue 'E30' '  +0.94%  '
Set-TextValue 'D31' '7.98'
Set-TextValue 'E31' '  -0.47%  '
Set-TextValue 'E32' '  -0.45%  '
Set-TextValue 'E33' '  +0.34%  '
Set-TextValue 'D34' '0.998'
Set-TextValue 'E34' '  +0.03%  '
Set-TextValue 'E35' '  -1.53%  '
Set-TextValue 'E36' '  +0.33%  '
Set-TextValue 'E37' '  +1.02%  '
Set-TextValue 'D38' '0.370'
Set-TextValue 'E38' '  +0.18%  '
Set-TextValue 'B39' 'EthereumClassic'
Set-TextValue 'C39' 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue 'D39' '18.79'
Set-TextValue 'E39' '  +0.26%  '
Set-TextValue 'D40' '151.76'
Set-TextValue 'E40' '  +0.93%  '
Set-TextValue 'B41' 'dogwifhat'
Set-TextValue 'C41' 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue 'D41' '2.78'
Set-TextValue 'E41' '  -0.02%  '
Set-TextValue 'D42' '1.76'
Set-TextValue 'E42' '  +1.03%  '
Set-TextValue 'E43' '  -0.03%  '
Set-TextValue 'B44' 'Aave'
Set-TextValue 'C44' 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue 'D44' '158.34'
Set-TextValue 'E44' '  +2.93%  '
Set-TextValue 'B45' 'BabyDogeCoin'
Set-TextValue 'C45' 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue 'D45' '0.0₆0303'
Set-TextValue 'E45' '  -3.00%  '
Set-TextValue 'D46' '15.41'
Set-TextValue 'E46' '  +1.10%  '
Set-TextValue 'E47' '  +0.84%  '
Set-TextValue 'D48' '20.83'
Set-TextValue 'E48' '  +1.27%  '
Set-TextValue 'D49' '0.615'
Set-TextValue 'E49' '  +1.16%  '
Set-TextValue 'D50' '0.0516'
Set-TextValue 'E50' '  +0.86%  '
Set-TextValue 'E51' '  +0.88%  '
